# Add the three new character styles referenced by the diff (GaNStyle,
# GaNParagraph, GaNLinks) and apply them to the runs that gained an
# <w:rPr><w:rStyle .../></w:rPr> in the commit.
#
# wdStyleTypeCharacter = 2
$d = $word.ActiveDocument

$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.NameAscii = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.NameAscii = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.NameAscii = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608   # RGB(0,0,128) -> w:color="000080"
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1     # wdUnderlineSingle

# --- Apply GaNStyle to every "Dates à utiliser ..." run (appears 4 times) ---
$datesTarget = "Dates à utiliser pour la Campagne 2022 Constellation de Persée: 16-25 janvier, 7-16 novembre, 6-15 décembre"
$rng = $d.Content
$rng.Find.ClearFormatting()
$datesStyled = 0
while ($rng.Find.Execute($datesTarget, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $datesStyled = $datesStyled + 1
    $rng.Collapse(0)
    if ($datesStyled -ge 10) { break }
}
Write-Host "GaNStyle applied to $datesStyled run(s)"

# --- Apply GaNParagraph to the "Vous allez participer ..." intro run ---
$paraTarget = "Vous allez participer à une campagne mondiale d’observation pour détecter les plus faibles étoiles visibles afin de mesurer la pollution lumineuse sur un site donné. Partout dans le monde, en localisant et en observant la Constellation de Persée dans le ciel nocturne et en la comparant aux cartes stellaires, les participants, apprendront comment l’éclairage, dans leur environnement local, influence la pollution lumineuse. Vos contributions à la base de données en ligne permettront de mesurer la qualité du ciel nocturne."
$rngParagraph = $d.Content
$rngParagraph.Find.ClearFormatting()
$paragraphFound = $rngParagraph.Find.Execute($paraTarget, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($paragraphFound) {
    $rngParagraph.Style = "GaNParagraph"
}
Write-Host "GaNParagraph applied: $paragraphFound"

# --- Apply GaNLinks to the "Les cartes figurant ..." credit run ---
$linksTarget = "Les cartes figurant dans ce document ont été établies par Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rngLinks = $d.Content
$rngLinks.Find.ClearFormatting()
$linksFound = $rngLinks.Find.Execute($linksTarget, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($linksFound) {
    $rngLinks.Style = "GaNLinks"
}
Write-Host "GaNLinks applied: $linksFound"
